$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I12").Value = "%"
$ws.Range("J12").Value = "Uninterpretable"
$ws.Range("I20").Value = "ba"
$ws.Range("J20").Value = "Appreciation"
$ws.Range("I30").Value = "sd"
$ws.Range("J30").Value = "Statement-non-opinion"
$ws.Range("I33").Value = "sd"
$ws.Range("J33").Value = "Statement-non-opinion"
$ws.Range("I49").Value = "b"
$ws.Range("J49").Value = "Acknowledge (Backchannel)"
$ws.Range("I51").Value = "aa"
$ws.Range("J51").Value = "Agree/Accept"
$ws.Range("I62").Value = "aa"
$ws.Range("J62").Value = "Agree/Accept"
$ws.Range("I71").Value = "b"
$ws.Range("J71").Value = "Acknowledge (Backchannel)"
$ws.Range("I75").Value = "aa"
$ws.Range("J75").Value = "Agree/Accept"
$ws.Range("I76").Value = "sd"
$ws.Range("J76").Value = "Statement-non-opinion"
$ws.Range("I77").Value = "%"
$ws.Range("J77").Value = "Uninterpretable"
$ws.Range("I91").Value = "sd"
$ws.Range("J91").Value = "Statement-non-opinion"
$ws.Range("I96").Value = "sv"
$ws.Range("J96").Value = "Statement-opinion"
$ws.Range("I97").Value = "aa"
$ws.Range("J97").Value = "Agree/Accept"
$ws.Range("I111").Value = "sv"
$ws.Range("J111").Value = "Statement-opinion"
$ws.Range("I113").Value = "sd"
$ws.Range("J113").Value = "Statement-non-opinion"
$ws.Range("I117").Value = "%"
$ws.Range("J117").Value = "Uninterpretable"
$ws.Range("I121").Value = "aa"
$ws.Range("J121").Value = "Agree/Accept"
$ws.Range("I123").Value = "sd"
$ws.Range("J123").Value = "Statement-non-opinion"
$ws.Range("I130").Value = "aa"
$ws.Range("J130").Value = "Agree/Accept"
$ws.Range("I132").Value = "sv"
$ws.Range("J132").Value = "Statement-opinion"
$ws.Range("I134").Value = "sv"
$ws.Range("J134").Value = "Statement-opinion"
$ws.Range("I137").Value = "aa"
$ws.Range("J137").Value = "Agree/Accept"
$ws.Range("I144").Value = "aa"
$ws.Range("J144").Value = "Agree/Accept"
$ws.Range("I149").Value = "%"
$ws.Range("J149").Value = "Uninterpretable"
$ws.Range("I159").Value = "%"
$ws.Range("J159").Value = "Uninterpretable"
$ws.Range("I160").Value = "b"
$ws.Range("J160").Value = "Acknowledge (Backchannel)"
$ws.Range("I166").Value = "aa"
$ws.Range("J166").Value = "Agree/Accept"
$ws.Range("I177").Value = "aa"
$ws.Range("J177").Value = "Agree/Accept"
